$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Certificate II in Split Air Conditioning and Heat Pump Systems
$ws.Range("C2").Value = "Air-Conditioning"
$ws.Range("M2").Value = "NSW/QLD"
$ws.Range("N2").Value = "Currently not accepting enrolments"
$ws.Range("R2").Value = ""
$ws.Rows.Item(2).RowHeight = 42.75

# Row 3 - Certificate IV in Air Conditioning and Refrigeration Servicing
$ws.Range("C3").Value = "Air-Conditioning"
$ws.Range("M3").Value = "NSW/QLD"
$ws.Range("N3").Value = "Currently not accepting enrolments"
$ws.Range("R3").Value = ""
$ws.Rows.Item(3).RowHeight = 42.75

# Row 4 - Advanced Diploma of Air Conditioning and Refrigeration Engineering
$ws.Range("C4").Value = "Air-Conditioning"
$ws.Range("M4").Value = "NSW/QLD"
$ws.Range("N4").Value = "Currently not accepting enrolments"
$ws.Range("R4").Value = ""
$ws.Rows.Item(4).RowHeight = 42.75

# Row 5 - Certificate III in Light Vehicle Mechanical Technology
$ws.Range("C5").Value = "Automotive"
$ws.Range("R5").Value = ""
$ws.Rows.Item(5).RowHeight = 42.75

# Row 6 - Certificate IV in Automotive Mechanical Diagnosis
$ws.Range("C6").Value = "Automotive"
$ws.Range("R6").Value = ""
$ws.Rows.Item(6).RowHeight = 42.75

# Row 7 - Certificate IV in Automotive Electrical Technology
$ws.Range("C7").Value = "Automotive"
$ws.Range("R7").Value = ""
$ws.Rows.Item(7).RowHeight = 42.75

# Row 8 - Diploma of Automotive Management
$ws.Range("C8").Value = "Automotive"
$ws.Range("R8").Value = ""
$ws.Rows.Item(8).RowHeight = 42.75

# Row 9 - Air-Conditioning package
$ws.Range("C9").Value = "Packages"
$ws.Range("M9").Value = "NSW/QLD"
$ws.Range("N9").Value = "Currently not accepting enrolments"
$ws.Range("R9").Value = ""
$ws.Rows.Item(9).RowHeight = 57

# Row 10 - Automotive package
$ws.Range("C10").Value = "Packages"
$ws.Range("R10").Value = ""
$ws.Rows.Item(10).RowHeight = 57

# Row 11 - Automotive package
$ws.Range("C11").Value = "Packages"
$ws.Range("R11").Value = ""
$ws.Rows.Item(11).RowHeight = 71.25

# Update selection to match the reviewed column (R2:R11)
$ws.Range("R2:R11").Select()
